$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.624.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'1.841.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'312.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.4249"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.07305"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'0.8725"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").Value = "'20.60"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'1.834.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "'5.329"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "'6.500"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").Value = "'0.06967"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'79.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.000008933"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "'15.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "'27.554.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'4.968"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'10.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").Value = "'2.053.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'1.981"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").Value = "'155.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'119.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("D29").Value = "'5.210"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").Value = "'1.860"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'0.08858"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "'0.7640"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "'2.956"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'4.495"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").Value = "'1.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").Value = "'0.05425"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "'2.818"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "'0.1660"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("D42").Value = "'0.5059"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'6.523"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.89%  "
$ws.Range("D44").Value = "'8.366"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.06544"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "'106.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.55%  "
$ws.Range("D47").Value = "'10.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'0.4626"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("D50").Value = "'1.634"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "'64.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.27%  "
